# Weekly data refresh: a new daily price record is inserted as row 74
# (sorted position within the existing date-ordered dataset), which pushes
# the previously-existing rows 74-166 down to 75-167.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 74; Excel shifts rows 74..166 down to 75..167
# and extends the used range/dimension to A1:T167 automatically.
$ws.Rows(74).Insert()

# Populate the newly inserted row 74 with the new record's data.
$ws.Range("A74").Value = 7
$ws.Range("B74").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C74").Value = "Ñuble"
$ws.Range("D74").Value = 45128
$ws.Range("E74").Value = 16
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100108
$ws.Range("H74").Value = "Tropicales y subtropicales"
$ws.Range("I74").Value = 100108002
$ws.Range("J74").Value = "Mango"
$ws.Range("K74").Value = "Sin especificar"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 40
$ws.Range("N74").Value = 9000
$ws.Range("O74").Value = 9000
$ws.Range("P74").Value = 9000
$ws.Range("Q74").Value = "`$/bandeja 4 kilos"
$ws.Range("R74").Value = "Brasil"
$ws.Range("S74").Value = 2250
$ws.Range("T74").Value = 4
